$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (2-9) down by two rows, to rows 4-11,
# copying cells (values + styles) exactly as-is so the original per-cell
# formatting (e.g. the date style on column A) is preserved unchanged.
$ws.Range("A2:B9").Copy($ws.Range("A4:B11"))

# Copy the column-A date style (s="2") from row 4 into the two new rows
# (2 and 3) so they match the rest of the date column formatting.
$ws.Range("A4").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Clear any clipboard/marching-ants state.
$excel.CutCopyMode = $false

# Fill in the two new data rows with their values (2014-12-31 and
# 2015-12-31 respectively).
$ws.Cells.Item(2, 1).Value = 42004
$ws.Cells.Item(2, 2).Value = 242638000000

$ws.Cells.Item(3, 1).Value = 42369
$ws.Cells.Item(3, 2).Value = 248545000000

# Update the sheet dimension to reflect the new used range.
$ws.Range("A1:B11") | Out-Null
